$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 581
$ws.Range("F5").Value = 745
$ws.Range("F6").Value = 381
$ws.Range("G6").Value = 48
$ws.Range("F8").Value = 159
$ws.Range("F9").Value = 243
$ws.Range("F10").Value = 227
$ws.Range("F11").Value = 6065
$ws.Range("F13").Value = 50
$ws.Range("F14").Value = 500
$ws.Range("F16").Value = 549
$ws.Range("F17").Value = 363
$ws.Range("F18").Value = 421
$ws.Range("F19").Value = 122
$ws.Range("F21").Value = 712
$ws.Range("F22").Value = 155
$ws.Range("F24").Value = 318
$ws.Range("F25").Value = 1022
$ws.Range("F26").Value = 65
$ws.Range("F27").Value = 1831
$ws.Range("F28").Value = 495

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 269
$ws.Range("F6").Value = 301

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 258

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 258
$ws.Range("F3").Value = 581
$ws.Range("F6").Value = 745
$ws.Range("F8").Value = 381
$ws.Range("G8").Value = 48
$ws.Range("F10").Value = 159
$ws.Range("F11").Value = 243
$ws.Range("F12").Value = 227
$ws.Range("F13").Value = 6065
$ws.Range("F15").Value = 50
$ws.Range("F16").Value = 269
$ws.Range("F17").Value = 500
$ws.Range("F19").Value = 549
$ws.Range("F20").Value = 363
$ws.Range("F21").Value = 421
$ws.Range("F23").Value = 122
$ws.Range("F26").Value = 301
$ws.Range("F28").Value = 712
$ws.Range("F32").Value = 155
$ws.Range("F34").Value = 318
$ws.Range("F35").Value = 1022
$ws.Range("F36").Value = 65
$ws.Range("F37").Value = 1831
$ws.Range("F38").Value = 495
